$d = $word.ActiveDocument

# Header row: "Group2" -> "Group0"
$d.Content.Find.Execute("Group2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Group0", 2)

# Remove the "Age - mean (sd)", "Quality of life - median (IQR)" and
# "Female - n (%)" summary rows from the first (and only) table.
$t = $d.Tables.Item(1)
$t.Rows.Item(2).Delete()
$t.Rows.Item(2).Delete()
$t.Rows.Item(2).Delete()
